$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price (D) cells whose new value would otherwise be auto-detected
# as a number by Excel stay plain text, exactly like the other text values in
# that column (e.g. "67.495.18").
$textPriceCells = @("D5","D6","D15","D19","D20","D21","D24","D25","D26","D28","D32","D34","D36","D38","D39","D40","D41","D42","D44","D45","D47","D48")
foreach ($addr in $textPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "67.495.18"
$ws.Range("E2").Value = "  -0.93%  "

# Row 3
$ws.Range("D3").Value = "3.225.18"
$ws.Range("E3").Value = "  -1.72%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").Value = "578.65"
$ws.Range("E5").Value = "  -1.55%  "

# Row 6
$ws.Range("D6").Value = "183.93"
$ws.Range("E6").Value = "  -1.54%  "

# Row 8
$ws.Range("E8").Value = "  +0.28%  "

# Row 9
$ws.Range("D9").Value = "3.223.28"
$ws.Range("E9").Value = "  -1.74%  "

# Row 10
$ws.Range("E10").Value = "  -3.17%  "

# Row 11
$ws.Range("E11").Value = "  -2.55%  "

# Row 12
$ws.Range("E12").Value = "  -1.66%  "

# Row 13
$ws.Range("D13").Value = "3.782.69"
$ws.Range("E13").Value = "  -1.74%  "

# Row 14
$ws.Range("E14").Value = "  +0.06%  "

# Row 15
$ws.Range("D15").Value = "27.74"
$ws.Range("E15").Value = "  -3.43%  "

# Row 16
$ws.Range("D16").Value = "67.526.35"
$ws.Range("E16").Value = "  -0.95%  "

# Row 17
$ws.Range("E17").Value = "  -2.20%  "

# Row 18
$ws.Range("D18").Value = "3.236.31"
$ws.Range("E18").Value = "  -1.29%  "

# Row 19
$ws.Range("D19").Value = "5.76"
$ws.Range("E19").Value = "  -2.03%  "

# Row 20
$ws.Range("D20").Value = "13.45"
$ws.Range("E20").Value = "  -1.73%  "

# Row 21
$ws.Range("D21").Value = "396.56"
$ws.Range("E21").Value = "  +3.48%  "

# Row 22
$ws.Range("E22").Value = "  -2.26%  "

# Row 23
$ws.Range("E23").Value = "  +0.10%  "

# Row 24
$ws.Range("D24").Value = "71.10"
$ws.Range("E24").Value = "  -0.65%  "

# Row 25
$ws.Range("D25").Value = "0.515"
$ws.Range("E25").Value = "  -0.45%  "

# Row 26
$ws.Range("D26").Value = "0.0000117"
$ws.Range("E26").Value = "  -3.17%  "

# Row 27
$ws.Range("E27").Value = "  -1.34%  "

# Row 28
$ws.Range("D28").Value = "9.55"
$ws.Range("E28").Value = "  -3.42%  "

# Row 29
$ws.Range("E29").Value = "  +0.19%  "

# Row 30
$ws.Range("E30").Value = "  -2.54%  "

# Row 31
$ws.Range("E31").Value = "  -4.86%  "

# Row 32
$ws.Range("D32").Value = "22.60"
$ws.Range("E32").Value = "  -1.55%  "

# Row 33
$ws.Range("E33").Value = "  -4.42%  "

# Row 34
$ws.Range("D34").Value = "0.998"
$ws.Range("E34").Value = "  +0.00%  "

# Row 35
$ws.Range("E35").Value = "  -3.11%  "

# Row 36
$ws.Range("D36").Value = "160.03"
$ws.Range("E36").Value = "  -1.33%  "

# Row 37
$ws.Range("E37").Value = "  -5.20%  "

# Row 38
$ws.Range("D38").Value = "1.89"
$ws.Range("E38").Value = "  +0.72%  "

# Row 39
$ws.Range("D39").Value = "26.35"
$ws.Range("E39").Value = "  -1.52%  "

# Row 40
$ws.Range("D40").Value = "0.803"
$ws.Range("E40").Value = "  -4.51%  "

# Row 41
$ws.Range("D41").Value = "4.54"
$ws.Range("E41").Value = "  -1.82%  "

# Row 42
$ws.Range("D42").Value = "6.50"
$ws.Range("E42").Value = "  -4.74%  "

# Row 43
$ws.Range("E43").Value = "  -6.02%  "

# Row 44
$ws.Range("D44").Value = "0.0683"
$ws.Range("E44").Value = "  -0.90%  "

# Row 45
$ws.Range("D45").Value = "40.60"
$ws.Range("E45").Value = "  -2.14%  "

# Row 46
$ws.Range("D46").Value = "2.592.34"
$ws.Range("E46").Value = "  -2.52%  "

# Row 47
$ws.Range("D47").Value = "24.46"
$ws.Range("E47").Value = "  -4.38%  "

# Row 48
$ws.Range("D48").Value = "332.72"
$ws.Range("E48").Value = "  -3.92%  "

# Row 49
$ws.Range("E49").Value = "  -3.01%  "

# Row 50
$ws.Range("E50").Value = "  -0.22%  "

# Row 51
$ws.Range("E51").Value = "  -2.11%  "
